$wb = $excel.ActiveWorkbook

$wsTodo = $wb.Worksheets.Item("TODO Before 0.0.1")
$wsLogs = $wb.Worksheets.Item("Logs")

# --- Sheet "TODO Before 0.0.1": add a new row 32 ---
$wsTodo.Activate()
$wsTodo.Range("B32").Value = "change cursor while casting some skill to some symbols (ua, lt)? "
$wsTodo.Range("B33").Select()

# --- Sheet "Logs": add two new rows (44, 45) ---
$wsLogs.Activate()

# Row 44 - copy formatting from row 43 (A: date style, B: wrap-text string style)
$wsLogs.Range("A43:B43").Copy()
$wsLogs.Range("A44:B44").PasteSpecial(-4122)
$wsLogs.Range("A44").Value = 45479
$wsLogs.Range("B44").Value = "oh…a lot of work with light, meshes gridmap, etc., no results only pain and not bad looking some kind of map"

# Row 45 - copy formatting from row 42 (taller, wrapped two-line row)
$wsLogs.Range("A42:B42").Copy()
$wsLogs.Range("A45:B45").PasteSpecial(-4122)
$wsLogs.Range("A45").Value = 45480
$wsLogs.Range("B45").Value = "gridMap removed - shitty staff. Work on random procedure generation of map. Have large progress with common Node3D and script for spawn tiles. Tiles could be spawn as rooms of different sizes, I can spawn long tunels of tiles with almost no collisions(low possibility). Next step apply rooms to tunels and add walls + roof"
$wsLogs.Rows.Item(45).RowHeight = 28.8

$wsLogs.Range("B46").Select()
